$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Summary")

# Enterprises density (per 1000 people): Micro / SMEs columns
$ws.Range("B11").NumberFormat = "@"
$ws.Range("B11").Value = "32.92"

$ws.Range("C11").NumberFormat = "@"
$ws.Range("C11").Value = "2.88"

# Employment (% of total): Micro / SMEs columns
$ws.Range("B12").NumberFormat = "@"
$ws.Range("B12").Value = "28.62"

$ws.Range("C12").NumberFormat = "@"
$ws.Range("C12").Value = "40.08"

# Enterprises (% of total): SMEs / MSMEs columns
$ws.Range("C36").NumberFormat = "@"
$ws.Range("C36").Value = "8.03"

$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = "99.72"
